# Mark attendance for step school - add 05-04-2025, 06-04-2025 (both "Off")
# and 07-04-2025 (actual attendance) columns, and bump Total Attendance
# for students who were Present on 07-04-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the three new header pairs (Status/Time) in R1:W1 ---------------
# Copy the style of the existing "01-04-2025 Time" header (Q1) onto the new
# header cells so they keep the same bold/centered/bordered formatting.
$ws.Range("Q1").Copy($ws.Range("R1:W1"))

$ws.Range("R1").Value = "05-04-2025 Status"
$ws.Range("S1").Value = "05-04-2025 Time"
$ws.Range("T1").Value = "06-04-2025 Status"
$ws.Range("U1").Value = "06-04-2025 Time"
$ws.Range("V1").Value = "07-04-2025 Status"
$ws.Range("W1").Value = "07-04-2025 Time"

# --- 2. Pre-create the blank "Time" cells for the two Off days --------------
# 05-04-2025 and 06-04-2025 are school holidays ("Off") for every student, so
# their Time columns (S and U) stay blank. Copy from a never-used, untouched
# column so the blank cells physically exist (matching a freshly attendance
# marked sheet) instead of simply being absent.
$ws.Range("Z1:Z23").Copy($ws.Range("S2:S24"))
$ws.Range("Z1:Z23").Copy($ws.Range("U2:U24"))

# --- 3. Per-student attendance data for the three new dates -----------------
# 05-04-2025 and 06-04-2025 are both marked "Off" (school holiday) for every
# student. 07-04-2025 carries the real attendance mark (P/A) and, for
# Present students, the check-in time.

$attendance = @{
    2  = @{ Status = "P"; Time = "09:42:38 AM" }
    3  = @{ Status = "P"; Time = "09:42:52 AM" }
    4  = @{ Status = "P"; Time = "09:42:56 AM" }
    5  = @{ Status = "P"; Time = "09:43:00 AM" }
    6  = @{ Status = "P"; Time = "09:43:10 AM" }
    7  = @{ Status = "A"; Time = "00:00:00" }
    8  = @{ Status = "P"; Time = "09:43:22 AM" }
    9  = @{ Status = "A"; Time = "00:00:00" }
    10 = @{ Status = "P"; Time = "09:43:26 AM" }
    11 = @{ Status = "A"; Time = "00:00:00" }
    12 = @{ Status = "A"; Time = "00:00:00" }
    13 = @{ Status = "A"; Time = "00:00:00" }
    14 = @{ Status = "P"; Time = "09:43:52 AM" }
    15 = @{ Status = "P"; Time = "09:43:55 AM" }
    16 = @{ Status = "P"; Time = "09:43:59 AM" }
    17 = @{ Status = "A"; Time = "00:00:00" }
    18 = @{ Status = "A"; Time = "00:00:00" }
    19 = @{ Status = "A"; Time = "00:00:00" }
    20 = @{ Status = "A"; Time = "00:00:00" }
    21 = @{ Status = "A"; Time = "00:00:00" }
    22 = @{ Status = "A"; Time = "00:00:00" }
    23 = @{ Status = "A"; Time = "00:00:00" }
    24 = @{ Status = "A"; Time = "00:00:00" }
}

foreach ($row in 2..24) {
    $info = $attendance[$row]

    # 05-04-2025 : Off day for everyone
    $ws.Cells.Item($row, 18).Value = "Off"   # R - Status

    # 06-04-2025 : Off day for everyone
    $ws.Cells.Item($row, 20).Value = "Off"   # T - Status

    # 07-04-2025 : real attendance mark
    $ws.Cells.Item($row, 22).Value = $info.Status  # V - Status
    $ws.Cells.Item($row, 23).Value = $info.Time     # W - Time

    # Bump the Total Attendance count for students present on 07-04-2025
    if ($info.Status -eq "P") {
        $ws.Cells.Item($row, 3).Value2 = $ws.Cells.Item($row, 3).Value2 + 1
    }
}
